# Append a new row to the (only) table in the document, recording the
# 14/01/2022 work-log entry: new map modes + start of province selection.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a new row at the end of the table. Word clones the tcPr/pPr
# formatting (borders, shading, spacing/justification) from the row
# that currently sits last, which matches the table's existing style.
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "14/01/2022"
$newRow.Cells.Item(2).Range.Text = "3 Hours 20 minutes"
$newRow.Cells.Item(3).Range.Text = "World Generation " + [char]0x2013 + " Objective 2"
$newRow.Cells.Item(4).Range.Text = "Added the map modes for showing elevation, temperature, rainfall and flora, as well as the storing and saving of these datapoints for provinces. Also began to work on implementing the single province viewer, but ran into issues with mesh hitboxes."
